# Dispatcher/Data/Config.xlsx re-organization
# - Move the "Work item" related settings (WorkItemWIID/WorkItemsPages/Status)
#   off of the Settings sheet into new Orchestrator Asset rows on the Assets sheet.
# - Keep only ProcessNames on the Settings sheet (shifted up into row 11).
# - Re-point which tab / cell is active & selected to match the new layout.

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")
$assets   = $wb.Worksheets.Item("Assets")

# --- Assets sheet : new rows for the work-item related asset values --------
# (written first, and in this specific order, so any newly-introduced text
# lands in the workbook's shared-string table in the same order it was
# authored in the original edit)
$assets.Range("B4").Value = "WIStatus"
$assets.Range("B2").Value = "WIType"
$assets.Range("B3").Value = "MaxPages"
$assets.Range("A4").Value = "Status"

$assets.Range("A2").Value = "WorkItemWIID"
$assets.Range("C2").Value = "AcmeTest"
$assets.Range("D2").Value = "Which wiid you want to extract"

$assets.Range("A3").Value = "WorkItemsPages"
$assets.Range("C3").Value = "AcmeTest"
$assets.Range("D3").Value = "The number of pages in work items section "

$assets.Range("C4").Value = "AcmeTest"

# --- Settings sheet ----------------------------------------------------------
# Row 11 used to hold "WorkItemWIID"; it now holds what used to be on row 13
# ("ProcessNames" / "excel;chrome"), and rows 12-13 are cleared out.
$settings.Range("A11").Value = "ProcessNames"
$settings.Range("B11").Value = "excel;chrome"
$settings.Range("C11").ClearContents()

$settings.Range("A12:C12").ClearContents()
$settings.Range("A13:C13").ClearContents()

# --- Active tab / selection -------------------------------------------------
# Settings is no longer the active tab; selection there moves to B3.
$settings.Activate() | Out-Null
$settings.Range("B3").Select() | Out-Null

# Assets becomes the active tab, with the selection sitting on A5.
$assets.Activate() | Out-Null
$assets.Range("A5").Select() | Out-Null
